# Apply the "new item added" edit to the DaySale report:
#  - Insert a new data row just above the existing "شاش 7 سم" row (row 25)
#    for a new item "جل فاتيكا اكياس" (Vatika gel sachets).
#  - Shift everything below it down by one row (Excel does this automatically
#    with Rows.Insert, including copying the per-row formatting/merges from
#    the row being pushed down).
#  - Update the running total (P column, last total row) to include the new
#    item's sale amount.
#  - Refresh the "printed at" timestamp in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 currently holds "شاش 7 سم" (and everything through row 32 is data,
# row 32 the subtotal, row 33 the footer). Insert a fresh row above it so the
# new item lands at row 25 and the rest of the table shifts down by one.
$ws.Rows("25:25").Insert()

# --- Fill in the new item row (row 25) ---
# Column A: sequence number (continues the existing 1..N numbering)
$ws.Cells.Item(25, 1).Value = 19
# Column C: item name
$ws.Cells.Item(25, 3).Value = "جل فاتيكا اكياس"
# Column H: current balance / quantity sold ("in:out" style counter)
$ws.Cells.Item(25, 8).Value = "47:0"
# Column L: order limit
$ws.Cells.Item(25, 12).Value = "0"
# Column N: unit price
$ws.Cells.Item(25, 14).Value = "2.50"
# Column P: sale price
$ws.Cells.Item(25, 16).Value = "5.0000"
# Column Q: number of transactions
$ws.Cells.Item(25, 17).Value = "2:0"

# --- Update the grand total row (now pushed down to row 33) ---
$ws.Cells.Item(33, 16).Value = 761.57

# --- Refresh the "printed at" timestamp in the footer (now row 34) ---
$ws.Cells.Item(34, 1).Value = "Saturday, 13 September, 2025 12:54 PM"
